$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: lengthen the underscore blanks in the "Copies served" line.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute(
    "Copies served by Dep. Clerk ___________ on the following date ___________ to:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Copies served by Dep. Clerk ___________________________ on the following date ____________________ to:",
    2)

# ---------------------------------------------------------------------------
# Edit 2: tweak the "Prosecutor's Office ..." service line.
#   a) one space is removed right after the first "PS"
#   b) a trailing semicolon is added at the very end of the line
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute(
    "Prosecutor’s Office: PS     OM",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Prosecutor’s Office: PS    OM",
    2)

$r3 = $d.Content
$r3.Find.Execute(
    "defendant.last_name}}: PS     OM     EM",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "defendant.last_name}}: PS     OM     EM;",
    2)

# ---------------------------------------------------------------------------
# Edit 3: insert a brand-new paragraph right after the "Prosecutor's Office"
# line containing the Community Control / County Jail conditional text.
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute(
    "Prosecutor’s Office",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($r4.Find.Found) {
    $officePara = $r4.Paragraphs(1)
    $anchorStart = $officePara.Range.Start
    $officePara.Range.InsertParagraphAfter()

    $allParas = $d.Paragraphs
    for ($i = 1; $i -le $allParas.Count; $i++) {
        $p = $allParas.Item($i)
        if ($p.Range.Start -eq $anchorStart) {
            $newPara = $allParas.Item($i + 1)
            $newPara.Range.Text = "{% if community_control.ordered is true or bond_conditions.monitoring is true %}Community Control: PS    EM;{% endif %}{% if jail_terms.ordered is true or apply_jtc == ‘Sentence’ %}County Jail: PS   EM;{% endif %}"
            break
        }
    }
}
